# Update crypto price/volume data per the latest scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "70.766.60"),
    @("E2", "  +1.50%  "),
    @("D3", "3.643.82"),
    @("E3", "  +6.31%  "),
    @("E4", "  +0.10%  "),
    @("D5", "581.64"),
    @("E5", "  -0.29%  "),
    @("D6", "176.27"),
    @("E6", "  -0.24%  "),
    @("D7", "3.637.11"),
    @("E7", "  +6.34%  "),
    @("D8", "0.617"),
    @("E8", "  +3.56%  "),
    @("D9", "0.999"),
    @("E9", "  -0.08%  "),
    @("E10", "  -1.53%  "),
    @("D11", "6.85"),
    @("E11", "  +24.17%  "),
    @("E12", "  +3.87%  "),
    @("D13", "48.69"),
    @("E13", "  -0.65%  "),
    @("E14", "  +2.58%  "),
    @("D15", "4.230.33"),
    @("E15", "  +6.45%  "),
    @("D16", "672.26"),
    @("E16", "  -2.63%  "),
    @("D17", "8.99"),
    @("E17", "  +4.30%  "),
    @("D18", "3.659.35"),
    @("E18", "  +6.87%  "),
    @("D19", "70.872.33"),
    @("E19", "  +1.62%  "),
    @("E20", "  +0.87%  "),
    @("D21", "17.82"),
    @("E21", "  +0.96%  "),
    @("D22", "11.50"),
    @("E22", "  +0.93%  "),
    @("D23", "0.935"),
    @("E23", "  +4.31%  "),
    @("D24", "17.21"),
    @("E24", "  +1.81%  "),
    @("D25", "101.06"),
    @("E25", "  +0.39%  "),
    @("E26", "  +0.44%  "),
    @("D27", "2.80"),
    @("E27", "  +5.73%  "),
    @("E29", "  -0.03%  "),
    @("D30", "34.89"),
    @("E30", "  +4.30%  "),
    @("E31", "  +4.11%  "),
    @("D32", "9.05"),
    @("E32", "  +3.69%  "),
    @("E33", "  -2.83%  "),
    @("D34", "7.54"),
    @("E34", "  +5.60%  "),
    @("E35", "  +6.72%  "),
    @("D36", "580.83"),
    @("E36", "  +1.25%  "),
    @("D37", "11.12"),
    @("E37", "  +1.13%  "),
    @("D38", "0.108"),
    @("E38", "  +4.39%  "),
    @("D39", "58.55"),
    @("E39", "  +0.43%  "),
    @("E40", "  +0.13%  "),
    @("D41", "3.598.47"),
    @("E41", "  +0.25%  "),
    @("E42", "  +8.82%  "),
    @("E43", "  +2.03%  "),
    @("D44", "0.348"),
    @("E44", "  +4.67%  "),
    @("E45", "  +1.54%  "),
    @("D46", "34.96"),
    @("E46", "  -0.36%  "),
    @("E47", "  +2.07%  "),
    @("E48", "  +8.41%  "),
    @("E49", "  +3.35%  "),
    @("D50", "134.90"),
    @("E50", "  +0.98%  "),
    @("D51", "2.95"),
    @("E51", "  +8.59%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $newValue = $u[1]
    $rng = $ws.Range($addr)

    # Force text storage so numeric-looking strings (e.g. "581.64",
    # "70.766.60") are written back as text, matching the original
    # inline-string cell type instead of being auto-coerced to a number.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}
